# Region IV-B_ALS-CLC.xlsx edit:
# Insert a new "INDEX (DO NOT MODIFY)" column at the very left of the sheet,
# shifting all existing columns one position to the right, uppercase all
# existing header labels (except the right-most "Status as of ..." column),
# and populate the new index column with a running index value (4 for this
# single data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column A; this shifts all existing content,
#    column widths, styles and data validation references one column right.
$ws.Columns.Item(1).Insert()

# 2. Give the new column A the same look (font/border/alignment) as the
#    header/data cells that used to live there (now in column B).
$ws.Range("B1:B2").Copy()
$ws.Range("A1:A2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 3. Set the new column's width and values. (ColumnWidth assignment in this
#    engine stores width + 5/6 character, so back the input off by 5/6 to
#    land exactly on a stored width of 23, matching the other clean widths.)
$ws.Columns.Item(1).ColumnWidth = 22.166666666666668
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"
$ws.Range("A2").Value = 4

# 4. Uppercase the pre-existing header labels (now shifted to columns B-X).
#    The very last header ("Status as of ...", now in column Y) is left as-is.
$ws.Range("B1").Value = "REGION"
$ws.Range("C1").Value = "DIVISION"
$ws.Range("D1").Value = "SCHOOL ID"
$ws.Range("E1").Value = "SCHOOL NAME"
$ws.Range("F1").Value = "MUNICIPALITY"
$ws.Range("G1").Value = "LEG DISTRICT"
$ws.Range("H1").Value = "NO. OF SITES"
$ws.Range("I1").Value = "SCOPE OF WORK"
$ws.Range("J1").Value = "TOTAL ALLOCATION"
$ws.Range("K1").Value = "CONTRACT AMOUNT"
$ws.Range("L1").Value = "STATUS"
$ws.Range("M1").Value = "PERCENTAGE OF COMPLETION"
$ws.Range("N1").Value = " TARGET COMPLETION DATE "
$ws.Range("O1").Value = "ACTUAL DATE OF COMPLETION"
$ws.Range("P1").Value = "PROJECT ID"
$ws.Range("Q1").Value = "CONTRACT ID"
$ws.Range("R1").Value = "ISSUANCE OF INVITATION TO BID"
$ws.Range("S1").Value = "PRE-SUBMISSION CONFERENCE"
$ws.Range("T1").Value = "BID OPENING"
$ws.Range("U1").Value = "ISSUANCE OF RESOLUTION TO AWARD"
$ws.Range("V1").Value = "ISSUANCE OF NOTICE TO PROCEED"
$ws.Range("W1").Value = "NAME OF CONTRACTOR"
$ws.Range("X1").Value = "OTHER REMARKS"
